$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add row 2: phone number 79174445 (as text), empty birthday, total_points 0
# Leading apostrophe forces text storage (matches source phone-number column
# being text), then reset the style so no explicit style index is left on
# the new cells (mirrors the target which has no s="..." on row 2).
$ws.Cells.Item(2, 1).Value = "'79174445"
$ws.Cells.Item(2, 1).Style = "Normal"

$ws.Cells.Item(2, 2).Value = "'"
$ws.Cells.Item(2, 2).Style = "Normal"

$ws.Cells.Item(2, 3).Value = 0
